$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D42", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = '29.242.94'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.862.30'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '0.7098'
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").Value = '237.77'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.08150'
$ws.Range("E8").Value = '  +9.20%  '
$ws.Range("D9").Value = '0.3038'
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").Value = '23.33'
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").Value = '0.08176'
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = '1.859.81'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '5.169'
$ws.Range("E13").Value = '  -1.38%  '
$ws.Range("D14").Value = '0.7088'
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").Value = '89.54'
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").Value = '29.248.07'
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").Value = '0.000007897'
$ws.Range("E17").Value = '  +3.51%  '
$ws.Range("D18").Value = '5.780'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '13.39'
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").Value = '236.91'
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '2.100.48'
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = '7.400'
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("D25").Value = '162.02'
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("D26").Value = '8.949'
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("D27").Value = '0.1456'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").Value = '18.08'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = '1.955'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("E30").Value = '  +2.08%  '
$ws.Range("D31").Value = '1.485'
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("D32").Value = '4.388'
$ws.Range("E32").Value = '  -3.61%  '
$ws.Range("D33").Value = '4.034'
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").Value = '0.05218'
$ws.Range("E34").Value = '  +0.91%  '
$ws.Range("D35").Value = '1.170'
$ws.Range("D36").Value = '0.7076'
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("D37").Value = '0.9992'
$ws.Range("E37").Value = '  -2.94%  '
$ws.Range("E38").Value = '  +0.50%  '
$ws.Range("D39").Value = '0.01854'
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("E40").Value = '  +1.63%  '
$ws.Range("D41").Value = '1.145.64'
$ws.Range("E41").Value = '  +6.61%  '
$ws.Range("D42").Value = '0.9242'
$ws.Range("E42").Value = '  -1.92%  '
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").Value = '5.868'
$ws.Range("E44").Value = '  -2.25%  '
$ws.Range("D45").Value = '70.23'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").Value = '0.9999'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '102.75'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("D48").Value = '1.773'
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("D49").Value = '1.999.87'
$ws.Range("E49").Value = '  +1.26%  '
$ws.Range("D50").Value = '9.222'
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("D51").Value = '6.955'
$ws.Range("E51").Value = '  -1.61%  '
